$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 72
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = "2024-07-02 17:17:00"
$ws.Range("D16").Value = 500
$ws.Range("E16").Value = "omelet "
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = "อาหาร"
$ws.Range("H16").Value = "expenses"
$ws.Range("I16").Value = "'"
$ws.Range("I16").Style = "Normal"
